$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.06427428475178958
$ws.Range("C2").Value = 0.7043999611326027
$ws.Range("D2").Value = 1.085782427793486
$ws.Range("E2").Value = 1.042008842473751
$ws.Range("F2").Value = 1.06570812435619

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1087012210651252
$ws.Range("C3").Value = 0.9180492773420234
$ws.Range("D3").Value = 1.469677934787732
$ws.Range("E3").Value = 1.212302740567607
$ws.Range("F3").Value = 1.235833313332441

# Row 4 (Q1)
$ws.Range("B4").Value = 0.09424136352415405
$ws.Range("C4").Value = 1.50554817053667
$ws.Range("D4").Value = 10.11276014084247
$ws.Range("E4").Value = 3.18005662541447
$ws.Range("F4").Value = 3.257157141059632
